$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet holds two side-by-side tables (分割前 / 分割後), each of which
# used to carry an extra "價值" (value) column, and the "分割後" table also
# had a trailing "差額" (difference) column. Drop them:
#
#   分割前 block: A:G -> A:F   (delete G, "價值")
#   分割後 block: I:P -> H:M   (delete the old O:P, "價值"/"差額", which
#                               after the first deletion above sit at N:O)

# Remove "價值" from the 分割前 (before-split) block.
$ws.Range("G1").EntireColumn.Delete()

# Remove "價值" and "差額" from the 分割後 (after-split) block — these used
# to be columns O:P, now shifted one to the left to N:O.
$ws.Range("N1:O1").EntireColumn.Delete()

# Re-apply the autofilter over the now-narrower A2:M34 range and keep the
# workbook's hidden _xlnm._FilterDatabase name in sync with it.
$ws.AutoFilterMode = $false
[void]$ws.Range("A2:M34").AutoFilter(1)

foreach ($n in $wb.Names) {
  if ($n.Name.EndsWith("_FilterDatabase")) {
    $n.RefersTo = "=" + $ws.Name + "!`$A`$2:`$M`$34"
  }
}

# Column A no longer needs to accommodate the wider original layout.
$ws.Columns("A:A").AutoFit()

# Match the reported final selection.
[void]$ws.Range("S10").Select()
